# Update the "Förändrad" (Changed) date column (C) for rows 2-51
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
